# Auto-generated cell updates applying refreshed market price data (H:N columns)
# across all 8 job sheets, per scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1150
$ws.Range("I4").Value = 1150
$ws.Range("K4").Value = 1150
$ws.Range("M4").Value = -1036
$ws.Range("H33").Value = 90.181816
$ws.Range("I33").Value = 54.6
$ws.Range("K33").Value = 54.6
$ws.Range("M33").Value = 174.4
$ws.Range("H39").Value = 92.59999999999999
$ws.Range("I39").Value = 95.75
$ws.Range("K39").Value = 287.25
$ws.Range("M39").Value = 8.75
$ws.Range("H86").Value = 1333.3334
$ws.Range("J86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 1333.3334
$ws.Range("J89").Value = 2000
$ws.Range("L89").Value = 10000
$ws.Range("N89").Value = -21232
$ws.Range("H98").Value = 1630
$ws.Range("I98").Value = 1714
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 1714
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = -216
$ws.Range("N98").Value = -3996
$ws.Range("H120").Value = 37748.75
$ws.Range("J120").Value = 37748.75
$ws.Range("L120").Value = 37748.75
$ws.Range("N120").Value = -47424.75
$ws.Range("H122").Value = 1630
$ws.Range("I122").Value = 1714
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5142
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2692
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 1221.7222
$ws.Range("I132").Value = 1240.6471
$ws.Range("K132").Value = 3721.9413
$ws.Range("M132").Value = -1191.9413
$ws.Range("H137").Value = 1850
$ws.Range("I137").Value = 1562.5
$ws.Range("K137").Value = 4687.5
$ws.Range("M137").Value = -2137.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3897.9443
$ws.Range("I32").Value = 2131.814
$ws.Range("J32").Value = 10801.909
$ws.Range("K32").Value = 2131.814
$ws.Range("L32").Value = 10801.909
$ws.Range("M32").Value = -1844.814
$ws.Range("N32").Value = -11375.909
$ws.Range("H45").Value = 1590.3334
$ws.Range("I45").Value = 1128.8334
$ws.Range("K45").Value = 1128.8334
$ws.Range("M45").Value = -751.8334
$ws.Range("H74").Value = 808.36664
$ws.Range("J74").Value = 1300
$ws.Range("L74").Value = 1300
$ws.Range("N74").Value = -3048
$ws.Range("H77").Value = 808.36664
$ws.Range("J77").Value = 1300
$ws.Range("L77").Value = 6500
$ws.Range("N77").Value = -15236
$ws.Range("H97").Value = 655.44446
$ws.Range("I97").Value = 612.375
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 612.375
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -116.375
$ws.Range("N97").Value = -1992
$ws.Range("H122").Value = 1526.25
$ws.Range("I122").Value = 1517.25
$ws.Range("J122").Value = 1553.25
$ws.Range("K122").Value = 4551.75
$ws.Range("L122").Value = 4659.75
$ws.Range("M122").Value = -2101.75
$ws.Range("N122").Value = -9559.75
$ws.Range("H132").Value = 1286.0322
$ws.Range("I132").Value = 1019.0417
$ws.Range("J132").Value = 2201.4285
$ws.Range("K132").Value = 3057.1251
$ws.Range("L132").Value = 6604.2855
$ws.Range("M132").Value = -527.1251000000002
$ws.Range("N132").Value = -11664.2855
$ws.Range("H135").Value = 42164.25
$ws.Range("J135").Value = 42164.25
$ws.Range("L135").Value = 42164.25
$ws.Range("N135").Value = -52304.25
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 111779.05
$ws.Range("J86").Value = 221784.9
$ws.Range("L86").Value = 221784.9
$ws.Range("N86").Value = -224030.9
$ws.Range("H89").Value = 111779.05
$ws.Range("J89").Value = 221784.9
$ws.Range("L89").Value = 1108924.5
$ws.Range("N89").Value = -1120156.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I7").Value = 131
$ws.Range("J7").Value = 63.75
$ws.Range("K7").Value = 131
$ws.Range("L7").Value = 63.75
$ws.Range("M7").Value = -18
$ws.Range("N7").Value = -289.75
$ws.Range("H31").Value = 3030.682
$ws.Range("I31").Value = 1956.9333
$ws.Range("J31").Value = 5331.5713
$ws.Range("K31").Value = 1956.9333
$ws.Range("L31").Value = 5331.5713
$ws.Range("M31").Value = -1661.9333
$ws.Range("N31").Value = -5921.5713
$ws.Range("H34").Value = 3030.682
$ws.Range("I34").Value = 1956.9333
$ws.Range("J34").Value = 5331.5713
$ws.Range("K34").Value = 1956.9333
$ws.Range("L34").Value = 5331.5713
$ws.Range("M34").Value = -1754.9333
$ws.Range("N34").Value = -5735.5713
$ws.Range("H74").Value = 29999.666
$ws.Range("J74").Value = 29999.666
$ws.Range("L74").Value = 29999.666
$ws.Range("N74").Value = -31747.666
$ws.Range("H77").Value = 29999.666
$ws.Range("J77").Value = 29999.666
$ws.Range("L77").Value = 89998.99800000001
$ws.Range("N77").Value = -98734.99800000001
$ws.Range("H122").Value = 1504.0358
$ws.Range("I122").Value = 1533.6111
$ws.Range("K122").Value = 4600.8333
$ws.Range("M122").Value = -2150.8333
$ws.Range("H134").Value = 948.1111
$ws.Range("I134").Value = 940.75
$ws.Range("J134").Value = 1007
$ws.Range("K134").Value = 2822.25
$ws.Range("L134").Value = 3021
$ws.Range("M134").Value = -287.25
$ws.Range("N134").Value = -8091

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 11470.228
$ws.Range("J4").Value = 1685033
$ws.Range("K4").Value = 34410.68399999999
$ws.Range("L4").Value = 5055099
$ws.Range("M4").Value = -34298.68399999999
$ws.Range("N4").Value = -5055323
$ws.Range("H9").Value = 262055.5
$ws.Range("I9").Value = 1000
$ws.Range("K9").Value = 3000
$ws.Range("M9").Value = -2776
$ws.Range("H37").Value = 99666
$ws.Range("J37").Value = 99666
$ws.Range("L37").Value = 298998
$ws.Range("N37").Value = -299222
$ws.Range("H50").Value = 83392424
$ws.Range("I50").Value = 140093.4
$ws.Range("J50").Value = 142858370
$ws.Range("K50").Value = 420280.2
$ws.Range("L50").Value = 428575110
$ws.Range("M50").Value = -419799.2
$ws.Range("N50").Value = -428576072
$ws.Range("H53").Value = 83392424
$ws.Range("I53").Value = 140093.4
$ws.Range("J53").Value = 142858370
$ws.Range("K53").Value = 420280.2
$ws.Range("L53").Value = 428575110
$ws.Range("M53").Value = -419799.2
$ws.Range("N53").Value = -428576072
$ws.Range("H107").Value = 1430.3125
$ws.Range("J107").Value = 1430.3125
$ws.Range("L107").Value = 4290.9375
$ws.Range("N107").Value = -8130.9375
$ws.Range("H131").Value = 812.3200000000001
$ws.Range("I131").Value = 511.42856
$ws.Range("J131").Value = 834.9677
$ws.Range("K131").Value = 1534.28568
$ws.Range("L131").Value = 2504.9031
$ws.Range("M131").Value = 3505.71432
$ws.Range("N131").Value = -12584.9031

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H102").Value = 2736.3125
$ws.Range("I102").Value = 2929.6667
$ws.Range("J102").Value = 2487.7144
$ws.Range("K102").Value = 2929.6667
$ws.Range("L102").Value = 2487.7144
$ws.Range("M102").Value = -1307.6667
$ws.Range("N102").Value = -5731.7144
$ws.Range("H132").Value = 1376083.4
$ws.Range("I132").Value = 1604124.6
$ws.Range("J132").Value = 7836.25
$ws.Range("K132").Value = 4812373.800000001
$ws.Range("L132").Value = 23508.75
$ws.Range("M132").Value = -4809843.800000001
$ws.Range("N132").Value = -28568.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1749.1765
$ws.Range("I132").Value = 1509.1578
$ws.Range("J132").Value = 2053.2
$ws.Range("K132").Value = 4527.4734
$ws.Range("L132").Value = 6159.599999999999
$ws.Range("M132").Value = -1997.4734
$ws.Range("N132").Value = -11219.6
$ws.Range("H136").Value = 3784
$ws.Range("I136").Value = 2120
$ws.Range("J136").Value = 5685.7144
$ws.Range("K136").Value = 6360
$ws.Range("L136").Value = 17057.1432
$ws.Range("M136").Value = -3810
$ws.Range("N136").Value = -22157.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 58833236
$ws.Range("J2").Value = 10312.6875
$ws.Range("L2").Value = 10312.6875
$ws.Range("N2").Value = -10536.6875
$ws.Range("H113").Value = 1283.4445
$ws.Range("I113").Value = 933.5
$ws.Range("K113").Value = 2800.5
$ws.Range("M113").Value = -630.5
$ws.Range("H126").Value = 7381.3213
$ws.Range("I126").Value = 6820
$ws.Range("J126").Value = 8248.817999999999
$ws.Range("K126").Value = 20460
$ws.Range("L126").Value = 24746.454
$ws.Range("M126").Value = -17990
$ws.Range("N126").Value = -29686.454
$ws.Range("H132").Value = 2517.0625
$ws.Range("I132").Value = 2231.25
$ws.Range("K132").Value = 6693.75
$ws.Range("M132").Value = -4163.75
$ws.Range("H136").Value = 15433442
$ws.Range("I136").Value = 20576994
$ws.Range("J136").Value = 2783
$ws.Range("K136").Value = 61730982
$ws.Range("L136").Value = 8349
$ws.Range("M136").Value = -61728432
$ws.Range("N136").Value = -13449

